$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the active cell selection to C9
$ws.Range("C9").Select()

# Rows 11 and 12: change B value from 1 to 2, and clear the fill (matches target style s=5: no fill, applyFill=1)
$ws.Range("B11").Value = 2
$ws.Range("B11").Interior.Pattern = -4142

$ws.Range("B12").Value = 2
$ws.Range("B12").Interior.Pattern = -4142

# Rows 18 and 19: change B value from 1 to 2, and clear the fill
$ws.Range("B18").Value = 2
$ws.Range("B18").Interior.Pattern = -4142

$ws.Range("B19").Value = 2
$ws.Range("B19").Interior.Pattern = -4142

# Rows 20 and 21: keep value at 1, but change fill color to theme color 9 (matches target style s=2)
$ws.Range("B20").Interior.ThemeColor = 9
$ws.Range("B21").Interior.ThemeColor = 9
